$wb = $excel.ActiveWorkbook

# Sheet 1: summ40133527 -> summ29536305
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ29536305"
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = "Intercept"
$arr[0,1] = [double]"8615.744493701804"
$arr[0,2] = [double]"2.125371173889367e-24"
$arr[1,0] = "Education[T.Primary/None]"
$arr[1,1] = [double]"-627.9508702055109"
$arr[1,2] = [double]"0.1146031708674052"
$arr[2,0] = "Education[T.Secondary]"
$arr[2,1] = [double]"-200.076703764995"
$arr[2,2] = [double]"0.5313126975616788"
$arr[3,0] = "Education[T.University]"
$arr[3,1] = [double]"47.42938155274899"
$arr[3,2] = [double]"0.7318457679435845"
$arr[4,0] = "Season[T.Spring]"
$arr[4,1] = [double]"-109.1390576782056"
$arr[4,2] = [double]"0.4520599776187583"
$arr[5,0] = "Season[T.Summer]"
$arr[5,1] = [double]"3.738630338087461"
$arr[5,2] = [double]"0.9814120544043627"
$arr[6,0] = "Season[T.Winter]"
$arr[6,1] = [double]"62.19715988411789"
$arr[6,2] = [double]"0.6696796369842908"
$arr[7,0] = "HHSize"
$arr[7,1] = [double]"-68.1947203293083"
$arr[7,2] = [double]"0.1595284879795326"
$arr[8,0] = "Sex"
$arr[8,1] = [double]"-1327.197989884121"
$arr[8,2] = [double]"5.400985602640942e-34"
$arr[9,0] = "Age"
$arr[9,1] = [double]"-24.80079946384163"
$arr[9,2] = [double]"1.053224548096937e-06"
$arr[10,0] = "DistSubcenter_res"
$arr[10,1] = [double]"71.39392503227299"
$arr[10,2] = [double]"0.06410379782752289"
$arr[11,0] = "DistCenter_res"
$arr[11,1] = [double]"639.2491693698491"
$arr[11,2] = [double]"3.51155502347245e-182"
$arr[12,0] = "UrbPopDensity_res"
$arr[12,1] = [double]"-0.03358217758785084"
$arr[12,2] = [double]"0.09154928862792104"
$arr[13,0] = "UrbBuildDensity_res"
$arr[13,1] = [double]"-1.212559955552506e-05"
$arr[13,2] = [double]"0.6991837992567951"
$arr[14,0] = "IntersecDensity_res"
$arr[14,1] = [double]"-17.39051690948888"
$arr[14,2] = [double]"0.02429329382402247"
$arr[15,0] = "street_length_res"
$arr[15,1] = [double]"-7.925438939866998"
$arr[15,2] = [double]"0.0409641780995532"
$arr[16,0] = "LU_Comm_res"
$arr[16,1] = [double]"49.51468668445159"
$arr[16,2] = [double]"0.9478840898998797"
$arr[17,0] = "LU_UrbFab_res"
$arr[17,1] = [double]"47.33970857695795"
$arr[17,2] = [double]"0.913637697027597"
$arr[18,0] = "bike_lane_share_res"
$arr[18,1] = [double]"-3096.930172037885"
$arr[18,2] = [double]"0.003218089376862115"
$ws.Range("A2:C20").Value2 = $arr

# Sheet 2: summ40863691 -> summ30097336
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ30097336"
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = "Intercept"
$arr[0,1] = [double]"8683.450129293491"
$arr[0,2] = [double]"1.494712333933632e-24"
$arr[1,0] = "Education[T.Primary/None]"
$arr[1,1] = [double]"-361.0469205321035"
$arr[1,2] = [double]"0.3694579389863428"
$arr[2,0] = "Education[T.Secondary]"
$arr[2,1] = [double]"-288.907895801201"
$arr[2,2] = [double]"0.3592922320258104"
$arr[3,0] = "Education[T.University]"
$arr[3,1] = [double]"-8.6236632131764"
$arr[3,2] = [double]"0.9499567508996769"
$arr[4,0] = "Season[T.Spring]"
$arr[4,1] = [double]"65.00511296858514"
$arr[4,2] = [double]"0.6518004557784085"
$arr[5,0] = "Season[T.Summer]"
$arr[5,1] = [double]"36.64959349619994"
$arr[5,2] = [double]"0.8182125075521204"
$arr[6,0] = "Season[T.Winter]"
$arr[6,1] = [double]"175.682269726226"
$arr[6,2] = [double]"0.2269741913867915"
$arr[7,0] = "HHSize"
$arr[7,1] = [double]"-36.87857478059482"
$arr[7,2] = [double]"0.4441036944477915"
$arr[8,0] = "Sex"
$arr[8,1] = [double]"-1326.748788622367"
$arr[8,2] = [double]"2.813262381429581e-34"
$arr[9,0] = "Age"
$arr[9,1] = [double]"-21.51330238991028"
$arr[9,2] = [double]"2.276799466458129e-05"
$arr[10,0] = "DistSubcenter_res"
$arr[10,1] = [double]"111.0853453692057"
$arr[10,2] = [double]"0.003650790910957478"
$arr[11,0] = "DistCenter_res"
$arr[11,1] = [double]"640.8480317526796"
$arr[11,2] = [double]"9.71275106119928e-184"
$arr[12,0] = "UrbPopDensity_res"
$arr[12,1] = [double]"-0.02364118484990374"
$arr[12,2] = [double]"0.2333494679570254"
$arr[13,0] = "UrbBuildDensity_res"
$arr[13,1] = [double]"-1.493919289669346e-05"
$arr[13,2] = [double]"0.6378690172079717"
$arr[14,0] = "IntersecDensity_res"
$arr[14,1] = [double]"-25.05698060951545"
$arr[14,2] = [double]"0.001087234849862768"
$arr[15,0] = "street_length_res"
$arr[15,1] = [double]"-10.30727582735412"
$arr[15,2] = [double]"0.008148984439376528"
$arr[16,0] = "LU_Comm_res"
$arr[16,1] = [double]"490.0479499620851"
$arr[16,2] = [double]"0.5185343251779556"
$arr[17,0] = "LU_UrbFab_res"
$arr[17,1] = [double]"13.21074523052641"
$arr[17,2] = [double]"0.9756982775259235"
$arr[18,0] = "bike_lane_share_res"
$arr[18,1] = [double]"-2810.399560546204"
$arr[18,2] = [double]"0.00777161612089739"
$ws.Range("A2:C20").Value2 = $arr

# Sheet 3: summ41665455 -> summ30709137
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ30709137"
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = "Intercept"
$arr[0,1] = [double]"8968.248488635581"
$arr[0,2] = [double]"1.674685146431559e-26"
$arr[1,0] = "Education[T.Primary/None]"
$arr[1,1] = [double]"-506.864155450565"
$arr[1,2] = [double]"0.2016820295918985"
$arr[2,0] = "Education[T.Secondary]"
$arr[2,1] = [double]"-321.4029962229469"
$arr[2,2] = [double]"0.3136201219590375"
$arr[3,0] = "Education[T.University]"
$arr[3,1] = [double]"91.67191068199483"
$arr[3,2] = [double]"0.5042429935450266"
$arr[4,0] = "Season[T.Spring]"
$arr[4,1] = [double]"-80.88557821912283"
$arr[4,2] = [double]"0.5736819833561959"
$arr[5,0] = "Season[T.Summer]"
$arr[5,1] = [double]"-127.0815750766818"
$arr[5,2] = [double]"0.423731613805903"
$arr[6,0] = "Season[T.Winter]"
$arr[6,1] = [double]"-15.35626114365047"
$arr[6,2] = [double]"0.9158693632328914"
$arr[7,0] = "HHSize"
$arr[7,1] = [double]"-27.73683303321791"
$arr[7,2] = [double]"0.5630500953541842"
$arr[8,0] = "Sex"
$arr[8,1] = [double]"-1336.54286720948"
$arr[8,2] = [double]"6.144529744108191e-35"
$arr[9,0] = "Age"
$arr[9,1] = [double]"-23.90877498741028"
$arr[9,2] = [double]"2.464206026912567e-06"
$arr[10,0] = "DistSubcenter_res"
$arr[10,1] = [double]"82.38871100626187"
$arr[10,2] = [double]"0.03148157986666132"
$arr[11,0] = "DistCenter_res"
$arr[11,1] = [double]"643.1149707245924"
$arr[11,2] = [double]"1.49109578817212e-186"
$arr[12,0] = "UrbPopDensity_res"
$arr[12,1] = [double]"-0.04572916504667483"
$arr[12,2] = [double]"0.01978011500810219"
$arr[13,0] = "UrbBuildDensity_res"
$arr[13,1] = [double]"1.214251156105532e-05"
$arr[13,2] = [double]"0.6942788866418231"
$arr[14,0] = "IntersecDensity_res"
$arr[14,1] = [double]"-23.1622070170839"
$arr[14,2] = [double]"0.002550445130775267"
$arr[15,0] = "street_length_res"
$arr[15,1] = [double]"-10.28082706144784"
$arr[15,2] = [double]"0.007473889901312926"
$arr[16,0] = "LU_Comm_res"
$arr[16,1] = [double]"-166.3315892467399"
$arr[16,2] = [double]"0.8246812875171872"
$arr[17,0] = "LU_UrbFab_res"
$arr[17,1] = [double]"142.5244515982102"
$arr[17,2] = [double]"0.7428306922012042"
$arr[18,0] = "bike_lane_share_res"
$arr[18,1] = [double]"-3902.292409667813"
$arr[18,2] = [double]"0.0001865640048103823"
$ws.Range("A2:C20").Value2 = $arr

# Sheet 4: summ42415121 -> summ31257904
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ31257904"
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = "Intercept"
$arr[0,1] = [double]"8170.794525245381"
$arr[0,2] = [double]"3.545758307557427e-22"
$arr[1,0] = "Education[T.Primary/None]"
$arr[1,1] = [double]"-582.0936115294874"
$arr[1,2] = [double]"0.1435800699512929"
$arr[2,0] = "Education[T.Secondary]"
$arr[2,1] = [double]"-530.68071538544"
$arr[2,2] = [double]"0.09270582982596108"
$arr[3,0] = "Education[T.University]"
$arr[3,1] = [double]"72.64486937572586"
$arr[3,2] = [double]"0.5978607649062506"
$arr[4,0] = "Season[T.Spring]"
$arr[4,1] = [double]"24.38494117664795"
$arr[4,2] = [double]"0.8658062722814852"
$arr[5,0] = "Season[T.Summer]"
$arr[5,1] = [double]"-42.32565190639681"
$arr[5,2] = [double]"0.7906435235233973"
$arr[6,0] = "Season[T.Winter]"
$arr[6,1] = [double]"25.4995566117486"
$arr[6,2] = [double]"0.8602954778886683"
$arr[7,0] = "HHSize"
$arr[7,1] = [double]"-43.65748111100959"
$arr[7,2] = [double]"0.3652878622098292"
$arr[8,0] = "Sex"
$arr[8,1] = [double]"-1232.745336521267"
$arr[8,2] = [double]"7.486509572490697e-30"
$arr[9,0] = "Age"
$arr[9,1] = [double]"-18.26650705067813"
$arr[9,2] = [double]"0.0003171406444963559"
$arr[10,0] = "DistSubcenter_res"
$arr[10,1] = [double]"68.29608896073398"
$arr[10,2] = [double]"0.0748089745159192"
$arr[11,0] = "DistCenter_res"
$arr[11,1] = [double]"645.6266396640322"
$arr[11,2] = [double]"3.5157020207708e-186"
$arr[12,0] = "UrbPopDensity_res"
$arr[12,1] = [double]"-0.03696969097620019"
$arr[12,2] = [double]"0.06130142953275322"
$arr[13,0] = "UrbBuildDensity_res"
$arr[13,1] = [double]"1.594562526122649e-05"
$arr[13,2] = [double]"0.6102196308983421"
$arr[14,0] = "IntersecDensity_res"
$arr[14,1] = [double]"-19.24164952423563"
$arr[14,2] = [double]"0.01201339191140978"
$arr[15,0] = "street_length_res"
$arr[15,1] = [double]"-7.362230142950185"
$arr[15,2] = [double]"0.05635777301594046"
$arr[16,0] = "LU_Comm_res"
$arr[16,1] = [double]"-238.7969616022515"
$arr[16,2] = [double]"0.7511419217297275"
$arr[17,0] = "LU_UrbFab_res"
$arr[17,1] = [double]"-163.2748620579959"
$arr[17,2] = [double]"0.7063696573732461"
$arr[18,0] = "bike_lane_share_res"
$arr[18,1] = [double]"-3173.027795697109"
$arr[18,2] = [double]"0.002487031143605859"
$ws.Range("A2:C20").Value2 = $arr

# Sheet 5: summ43173716 -> summ31825661
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ31825661"
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = "Intercept"
$arr[0,1] = [double]"7930.333368268832"
$arr[0,2] = [double]"8.678342730102076e-21"
$arr[1,0] = "Education[T.Primary/None]"
$arr[1,1] = [double]"-622.8480726436673"
$arr[1,2] = [double]"0.122234540318375"
$arr[2,0] = "Education[T.Secondary]"
$arr[2,1] = [double]"-199.7332017537057"
$arr[2,2] = [double]"0.5271714246335852"
$arr[3,0] = "Education[T.University]"
$arr[3,1] = [double]"78.61615814687555"
$arr[3,2] = [double]"0.568418315348113"
$arr[4,0] = "Season[T.Spring]"
$arr[4,1] = [double]"88.14161196131508"
$arr[4,2] = [double]"0.5427425818212535"
$arr[5,0] = "Season[T.Summer]"
$arr[5,1] = [double]"99.61695807996611"
$arr[5,2] = [double]"0.5352035336484944"
$arr[6,0] = "Season[T.Winter]"
$arr[6,1] = [double]"100.1771161239493"
$arr[6,2] = [double]"0.4932107873195206"
$arr[7,0] = "HHSize"
$arr[7,1] = [double]"-86.74346048844032"
$arr[7,2] = [double]"0.07428510025348169"
$arr[8,0] = "Sex"
$arr[8,1] = [double]"-1175.448733255416"
$arr[8,2] = [double]"5.497121269803767e-27"
$arr[9,0] = "Age"
$arr[9,1] = [double]"-25.33846607974916"
$arr[9,2] = [double]"6.431622857872734e-07"
$arr[10,0] = "DistSubcenter_res"
$arr[10,1] = [double]"88.31269993166097"
$arr[10,2] = [double]"0.02082250493144727"
$arr[11,0] = "DistCenter_res"
$arr[11,1] = [double]"623.9486942608974"
$arr[11,2] = [double]"6.389223592439651e-174"
$arr[12,0] = "UrbPopDensity_res"
$arr[12,1] = [double]"-0.04198638393091304"
$arr[12,2] = [double]"0.03556897577224163"
$arr[13,0] = "UrbBuildDensity_res"
$arr[13,1] = [double]"-1.36660633916929e-06"
$arr[13,2] = [double]"0.9657930527373269"
$arr[14,0] = "IntersecDensity_res"
$arr[14,1] = [double]"-15.53950940386902"
$arr[14,2] = [double]"0.045368422442627"
$arr[15,0] = "street_length_res"
$arr[15,1] = [double]"-3.359692340370085"
$arr[15,2] = [double]"0.38831076216635"
$arr[16,0] = "LU_Comm_res"
$arr[16,1] = [double]"169.6235410483932"
$arr[16,2] = [double]"0.8224053455582045"
$arr[17,0] = "LU_UrbFab_res"
$arr[17,1] = [double]"275.5683300967787"
$arr[17,2] = [double]"0.5273550376144482"
$arr[18,0] = "bike_lane_share_res"
$arr[18,1] = [double]"-4292.955125302722"
$arr[18,2] = [double]"5.031594662729359e-05"
$ws.Range("A2:C20").Value2 = $arr

# Sheet 6: summ43948294 -> summ32379183
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ32379183"
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = "Intercept"
$arr[0,1] = [double]"8507.955097993552"
$arr[0,2] = [double]"7.52316808338126e-24"
$arr[1,0] = "Education[T.Primary/None]"
$arr[1,1] = [double]"-294.3572243150239"
$arr[1,2] = [double]"0.4685831906985024"
$arr[2,0] = "Education[T.Secondary]"
$arr[2,1] = [double]"110.8091442479201"
$arr[2,2] = [double]"0.72871326415625"
$arr[3,0] = "Education[T.University]"
$arr[3,1] = [double]"80.76584033529683"
$arr[3,2] = [double]"0.5572240858295551"
$arr[4,0] = "Season[T.Spring]"
$arr[4,1] = [double]"-24.76446540845099"
$arr[4,2] = [double]"0.8646614982063865"
$arr[5,0] = "Season[T.Summer]"
$arr[5,1] = [double]"44.24561634992338"
$arr[5,2] = [double]"0.7837165451302734"
$arr[6,0] = "Season[T.Winter]"
$arr[6,1] = [double]"98.62399768402368"
$arr[6,2] = [double]"0.5000735568973531"
$arr[7,0] = "HHSize"
$arr[7,1] = [double]"-67.07670635238983"
$arr[7,2] = [double]"0.1673614941704978"
$arr[8,0] = "Sex"
$arr[8,1] = [double]"-1362.077324134767"
$arr[8,2] = [double]"1.663176870529806e-35"
$arr[9,0] = "Age"
$arr[9,1] = [double]"-23.17110644946975"
$arr[9,2] = [double]"6.149920708499743e-06"
$arr[10,0] = "DistSubcenter_res"
$arr[10,1] = [double]"71.94847895065038"
$arr[10,2] = [double]"0.06015144616034981"
$arr[11,0] = "DistCenter_res"
$arr[11,1] = [double]"660.9432768598417"
$arr[11,2] = [double]"6.657768168547409e-193"
$arr[12,0] = "UrbPopDensity_res"
$arr[12,1] = [double]"-0.02010990762295771"
$arr[12,2] = [double]"0.3123271781451656"
$arr[13,0] = "UrbBuildDensity_res"
$arr[13,1] = [double]"-1.340589141347778e-05"
$arr[13,2] = [double]"0.6693341458523896"
$arr[14,0] = "IntersecDensity_res"
$arr[14,1] = [double]"-18.78896236800118"
$arr[14,2] = [double]"0.01510638445552061"
$arr[15,0] = "street_length_res"
$arr[15,1] = [double]"-10.18044997489479"
$arr[15,2] = [double]"0.008213809016200264"
$arr[16,0] = "LU_Comm_res"
$arr[16,1] = [double]"566.5146111194153"
$arr[16,2] = [double]"0.456311898391246"
$arr[17,0] = "LU_UrbFab_res"
$arr[17,1] = [double]"-81.90674117797164"
$arr[17,2] = [double]"0.8513318580657605"
$arr[18,0] = "bike_lane_share_res"
$arr[18,1] = [double]"-3072.95390480251"
$arr[18,2] = [double]"0.003844506431402041"
$ws.Range("A2:C20").Value2 = $arr

# Sheet 7: summ44722097 -> summ32943600
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ32943600"
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = "Intercept"
$arr[0,1] = [double]"9138.918585827745"
$arr[0,2] = [double]"4.258839165617422e-27"
$arr[1,0] = "Education[T.Primary/None]"
$arr[1,1] = [double]"-534.5967780269036"
$arr[1,2] = [double]"0.1837573430245569"
$arr[2,0] = "Education[T.Secondary]"
$arr[2,1] = [double]"1.141722250464682"
$arr[2,2] = [double]"0.9971086225912389"
$arr[3,0] = "Education[T.University]"
$arr[3,1] = [double]"159.9163962509835"
$arr[3,2] = [double]"0.242179846682603"
$arr[4,0] = "Season[T.Spring]"
$arr[4,1] = [double]"32.79972056908807"
$arr[4,2] = [double]"0.8200366327575974"
$arr[5,0] = "Season[T.Summer]"
$arr[5,1] = [double]"-12.7968498485271"
$arr[5,2] = [double]"0.9360951871545642"
$arr[6,0] = "Season[T.Winter]"
$arr[6,1] = [double]"76.36512615944906"
$arr[6,2] = [double]"0.5993021044981185"
$arr[7,0] = "HHSize"
$arr[7,1] = [double]"-28.62723747723759"
$arr[7,2] = [double]"0.5536955295823629"
$arr[8,0] = "Sex"
$arr[8,1] = [double]"-1296.374785534196"
$arr[8,2] = [double]"8.515783598503956e-33"
$arr[9,0] = "Age"
$arr[9,1] = [double]"-22.8401604211209"
$arr[9,2] = [double]"7.105020075784274e-06"
$arr[10,0] = "DistSubcenter_res"
$arr[10,1] = [double]"97.1458683229954"
$arr[10,2] = [double]"0.0105628527378144"
$arr[11,0] = "DistCenter_res"
$arr[11,1] = [double]"649.0622225980387"
$arr[11,2] = [double]"2.078005725358033e-189"
$arr[12,0] = "UrbPopDensity_res"
$arr[12,1] = [double]"-0.03297581494240699"
$arr[12,2] = [double]"0.09944757902038123"
$arr[13,0] = "UrbBuildDensity_res"
$arr[13,1] = [double]"-1.691913810796843e-06"
$arr[13,2] = [double]"0.9578575082341856"
$arr[14,0] = "IntersecDensity_res"
$arr[14,1] = [double]"-25.02136635268098"
$arr[14,2] = [double]"0.001172103057971376"
$arr[15,0] = "street_length_res"
$arr[15,1] = [double]"-14.72347382809482"
$arr[15,2] = [double]"0.0001484891729436207"
$arr[16,0] = "LU_Comm_res"
$arr[16,1] = [double]"656.7047212703037"
$arr[16,2] = [double]"0.3868117020891888"
$arr[17,0] = "LU_UrbFab_res"
$arr[17,1] = [double]"193.6717848361453"
$arr[17,2] = [double]"0.6540225709000012"
$arr[18,0] = "bike_lane_share_res"
$arr[18,1] = [double]"-3648.730569636733"
$arr[18,2] = [double]"0.0004742796676337329"
$ws.Range("A2:C20").Value2 = $arr

# Sheet 8: summ45492432 -> summ33578192
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ33578192"
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = "Intercept"
$arr[0,1] = [double]"8708.601249987123"
$arr[0,2] = [double]"2.652180676831744e-25"
$arr[1,0] = "Education[T.Primary/None]"
$arr[1,1] = [double]"-519.0243940949634"
$arr[1,2] = [double]"0.2118574930597281"
$arr[2,0] = "Education[T.Secondary]"
$arr[2,1] = [double]"-348.1569770141797"
$arr[2,2] = [double]"0.2673132880918865"
$arr[3,0] = "Education[T.University]"
$arr[3,1] = [double]"12.12759944165036"
$arr[3,2] = [double]"0.9294714909440537"
$arr[4,0] = "Season[T.Spring]"
$arr[4,1] = [double]"-12.5329637164662"
$arr[4,2] = [double]"0.9307644482269756"
$arr[5,0] = "Season[T.Summer]"
$arr[5,1] = [double]"-56.36768098263829"
$arr[5,2] = [double]"0.7244857097116801"
$arr[6,0] = "Season[T.Winter]"
$arr[6,1] = [double]"85.795389171712"
$arr[6,2] = [double]"0.5539906064198218"
$arr[7,0] = "HHSize"
$arr[7,1] = [double]"-49.71154461179242"
$arr[7,2] = [double]"0.303297866686833"
$arr[8,0] = "Sex"
$arr[8,1] = [double]"-1371.318258698805"
$arr[8,2] = [double]"1.597798198217971e-36"
$arr[9,0] = "Age"
$arr[9,1] = [double]"-23.49188917268371"
$arr[9,2] = [double]"3.536144797654455e-06"
$arr[10,0] = "DistSubcenter_res"
$arr[10,1] = [double]"114.5659724817248"
$arr[10,2] = [double]"0.002599582320034618"
$arr[11,0] = "DistCenter_res"
$arr[11,1] = [double]"650.5858755907765"
$arr[11,2] = [double]"2.476534438950354e-191"
$arr[12,0] = "UrbPopDensity_res"
$arr[12,1] = [double]"-0.02660871112222504"
$arr[12,2] = [double]"0.1786081552041299"
$arr[13,0] = "UrbBuildDensity_res"
$arr[13,1] = [double]"4.323970838268913e-07"
$arr[13,2] = [double]"0.9890255624667488"
$arr[14,0] = "IntersecDensity_res"
$arr[14,1] = [double]"-20.99301766122095"
$arr[14,2] = [double]"0.006160929268875907"
$arr[15,0] = "street_length_res"
$arr[15,1] = [double]"-10.80593695622412"
$arr[15,2] = [double]"0.004613487574517559"
$arr[16,0] = "LU_Comm_res"
$arr[16,1] = [double]"812.8998284690506"
$arr[16,2] = [double]"0.2793319358829496"
$arr[17,0] = "LU_UrbFab_res"
$arr[17,1] = [double]"-79.54354832027843"
$arr[17,2] = [double]"0.8539527608728616"
$arr[18,0] = "bike_lane_share_res"
$arr[18,1] = [double]"-3684.289178173742"
$arr[18,2] = [double]"0.0003952964141797342"
$ws.Range("A2:C20").Value2 = $arr

# Sheet 9: summ46281449 -> summ34152988
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ34152988"
$arr = New-Object 'object[,]' 19,3
$arr[0,0] = "Intercept"
$arr[0,1] = [double]"8707.068235132114"
$arr[0,2] = [double]"8.830759036534143e-25"
$arr[1,0] = "Education[T.Primary/None]"
$arr[1,1] = [double]"-884.6538810508962"
$arr[1,2] = [double]"0.03008885505549092"
$arr[2,0] = "Education[T.Secondary]"
$arr[2,1] = [double]"-231.7861391465376"
$arr[2,2] = [double]"0.4657489048792317"
$arr[3,0] = "Education[T.University]"
$arr[3,1] = [double]"-95.99451913846102"
$arr[3,2] = [double]"0.4864130416552213"
$arr[4,0] = "Season[T.Spring]"
$arr[4,1] = [double]"84.67090904542752"
$arr[4,2] = [double]"0.5583568221320614"
$arr[5,0] = "Season[T.Summer]"
$arr[5,1] = [double]"-20.59143692770172"
$arr[5,2] = [double]"0.897085497654409"
$arr[6,0] = "Season[T.Winter]"
$arr[6,1] = [double]"127.2576517434841"
$arr[6,2] = [double]"0.3819627577687728"
$arr[7,0] = "HHSize"
$arr[7,1] = [double]"-38.98269811033661"
$arr[7,2] = [double]"0.4183961936051797"
$arr[8,0] = "Sex"
$arr[8,1] = [double]"-1298.75276881628"
$arr[8,2] = [double]"8.672446053168025e-33"
$arr[9,0] = "Age"
$arr[9,1] = [double]"-26.46467419403063"
$arr[9,2] = [double]"1.809481878272657e-07"
$arr[10,0] = "DistSubcenter_res"
$arr[10,1] = [double]"118.0820078626256"
$arr[10,2] = [double]"0.002026765468874975"
$arr[11,0] = "DistCenter_res"
$arr[11,1] = [double]"637.1032044673626"
$arr[11,2] = [double]"1.502105649949431e-179"
$arr[12,0] = "UrbPopDensity_res"
$arr[12,1] = [double]"-0.02365867256164274"
$arr[12,2] = [double]"0.2352339631758779"
$arr[13,0] = "UrbBuildDensity_res"
$arr[13,1] = [double]"-6.477949339452169e-06"
$arr[13,2] = [double]"0.8389460336670985"
$arr[14,0] = "IntersecDensity_res"
$arr[14,1] = [double]"-18.54982324576741"
$arr[14,2] = [double]"0.01656526853990092"
$arr[15,0] = "street_length_res"
$arr[15,1] = [double]"-9.493279895067463"
$arr[15,2] = [double]"0.01406646622380662"
$arr[16,0] = "LU_Comm_res"
$arr[16,1] = [double]"639.4451043530885"
$arr[16,2] = [double]"0.3984192151255548"
$arr[17,0] = "LU_UrbFab_res"
$arr[17,1] = [double]"77.31845248779462"
$arr[17,2] = [double]"0.8587303746763724"
$arr[18,0] = "bike_lane_share_res"
$arr[18,1] = [double]"-4509.596828803138"
$arr[18,2] = [double]"1.67369186614251e-05"
$ws.Range("A2:C20").Value2 = $arr
